$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format while writing so that numeric-looking
# strings ("1.00", "368.78", ...) are not silently coerced into Double values,
# matching the original inlineStr/text cells. ClearFormats() afterwards restores
# the default (unstyled) cell format so no stray style index is left behind.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '50.861.88'
$ws.Range("D3").Value = '2.901.15'
$ws.Range("D4").Value = '1.00'
$ws.Range("D5").Value = '368.78'
$ws.Range("D6").Value = '102.01'
$ws.Range("D7").Value = '0.539'
$ws.Range("D9").Value = '0.580'
$ws.Range("D10").Value = '36.62'
$ws.Range("D12").Value = '0.0831'
$ws.Range("D13").Value = '18.21'
$ws.Range("D14").Value = '3.353.56'
$ws.Range("D15").Value = '7.35'
$ws.Range("D16").Value = '2.895.87'
$ws.Range("D17").Value = '0.918'
$ws.Range("D18").Value = '50.826.60'
$ws.Range("D20").Value = '7.15'
$ws.Range("D21").Value = '12.80'
$ws.Range("D22").Value = '0.0₃0938'
$ws.Range("D23").Value = '67.76'
$ws.Range("D24").Value = '257.49'
$ws.Range("D25").Value = '2.65'
$ws.Range("D29").Value = '25.48'
$ws.Range("D32").Value = '6.21'
$ws.Range("D33").Value = '9.83'
$ws.Range("D35").Value = '51.27'
$ws.Range("D36").Value = '33.89'
$ws.Range("D38").Value = '0.0417'
$ws.Range("D40").Value = '16.93'
$ws.Range("D41").Value = '2.56'
$ws.Range("D43").Value = '0.111'
$ws.Range("D44").Value = '118.92'
$ws.Range("D45").Value = '21.76'
$ws.Range("D46").Value = '2.09'
$ws.Range("D47").Value = '2.014.73'
$ws.Range("D49").Value = '3.11'
$ws.Range("D50").Value = '3.189.29'

$dRange.ClearFormats()

$ws.Range("E2").Formula = '  -1.90%  '
$ws.Range("E3").Formula = '  -2.00%  '
$ws.Range("E4").Formula = '  -0.12%  '
$ws.Range("E5").Formula = '  +4.72%  '
$ws.Range("E6").Formula = '  -4.38%  '
$ws.Range("E7").Formula = '  -3.07%  '
$ws.Range("E8").Formula = '  -0.11%  '
$ws.Range("E9").Formula = '  -4.70%  '
$ws.Range("E10").Formula = '  -3.82%  '
$ws.Range("E11").Formula = '  +0.36%  '
$ws.Range("E12").Formula = '  -2.25%  '
$ws.Range("E13").Formula = '  -4.34%  '
$ws.Range("E14").Formula = '  -2.03%  '
$ws.Range("E15").Formula = '  -3.15%  '
$ws.Range("E16").Formula = '  -2.56%  '
$ws.Range("E17").Formula = '  -7.32%  '
$ws.Range("E18").Formula = '  -1.97%  '
$ws.Range("E19").Formula = '  -6.28%  '
$ws.Range("E20").Formula = '  -3.74%  '
$ws.Range("E21").Formula = '  -4.93%  '
$ws.Range("E22").Formula = '  -3.21%  '
$ws.Range("E23").Formula = '  -2.12%  '
$ws.Range("E24").Formula = '  -1.91%  '
$ws.Range("E25").Formula = '  -2.44%  '
$ws.Range("E26").Formula = '  -2.58%  '
$ws.Range("E28").Formula = '  -4.86%  '
$ws.Range("E29").Formula = '  -4.30%  '
$ws.Range("E30").Formula = '  -3.71%  '
$ws.Range("E31").Formula = '  -4.99%  '
$ws.Range("E32").Formula = '  +2.34%  '
$ws.Range("E33").Formula = '  -4.28%  '
$ws.Range("E34").Formula = '  -2.34%  '
$ws.Range("E35").Formula = '  +1.39%  '
$ws.Range("E36").Formula = '  -5.74%  '
$ws.Range("E37").Formula = '  +0.63%  '
$ws.Range("E38").Formula = '  -3.69%  '
$ws.Range("E39").Formula = '  -6.83%  '
$ws.Range("E40").Formula = '  -4.63%  '
$ws.Range("E41").Formula = '  -4.49%  '
$ws.Range("E42").Formula = '  -6.30%  '
$ws.Range("E43").Formula = '  -3.69%  '
$ws.Range("E44").Formula = '  -3.72%  '
$ws.Range("E45").Formula = '  -2.92%  '
$ws.Range("E46").Formula = '  -1.06%  '
$ws.Range("E47").Formula = '  -4.57%  '
$ws.Range("E48").Formula = '  -0.12%  '
$ws.Range("E49").Formula = '  -6.45%  '
$ws.Range("E50").Formula = '  -1.60%  '
$ws.Range("E51").Formula = '  -0.72%  '
